$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$table.Columns.Item(1).Width = 953 / 20
$table.Columns.Item(2).Width = 953 / 20
$table.Columns.Item(3).Width = 881 / 20
$table.Columns.Item(4).Width = 863 / 20
$table.Columns.Item(5).Width = 863 / 20
$table.Columns.Item(6).Width = 863 / 20
$table.Columns.Item(7).Width = 863 / 20
$table.Columns.Item(8).Width = 863 / 20
$table.Columns.Item(9).Width = 863 / 20
$table.Columns.Item(10).Width = 863 / 20

$table.Cell(1, 1).Range.Text = "Numero"
$table.Cell(1, 2).Range.Text = "Nombre"
$table.Cell(1, 3).Range.Text = "Carnet"
